$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Row 39 gets finalized: Status -> DONE, Finalized date filled in,
# and a new Fee cell is added.
# ------------------------------------------------------------------
$ws.Cells.Item(39, 8).Value = "DONE"
$ws.Cells.Item(39, 9).Value = 42859.875844907408
$ws.Cells.Item(39, 10).Value = "0.00510000 ETC (0.15%)"

# ------------------------------------------------------------------
# New row 40: the next (still in-progress) transaction.
# ------------------------------------------------------------------

# A40 - Data (date/time), same numeric format as A39 (style index 6:
# m/d/yy h:mm + wrap text).
$ws.Cells.Item(40, 1).NumberFormat = "m/d/yy h:mm"
$ws.Cells.Item(40, 1).WrapText = $true
$ws.Cells.Item(40, 1).Value = 42860.441886574074

# B40 - Action(Buy/Sell): reuse the existing rich-text "Sell" shared
# string (red text) by copying an existing Sell cell's value/format.
$ws.Range("B2").Copy()
$ws.Range("B40").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# C40 - Currency
$ws.Cells.Item(40, 3).Value = "        ETC"

# D40 - Transaction code: this is a purely numeric-looking string that
# must stay stored as literal text (matches the existing shared string
# table layout). Build it as a formula first so Excel treats it as
# text, then convert the formula to a plain value via copy / paste-
# special so the underlying cell type stays "text" instead of being
# coerced to a number.
$ws.Cells.Item(40, 4).Formula = '="               7.50999943"&CHAR(10)'
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Cells.Item(40, 4).WrapText = $true

# E40 - Transaction value
$ws.Cells.Item(40, 5).Value = "          7.7USDT"

# F40 - Transaction amount
$ws.Cells.Item(40, 6).Value = "        3.4ETC"

# G40 - Current value (USDT) / transaction code reference
$ws.Cells.Item(40, 7).Value = " ETC/USDT0000001"

# H40 - Status
$ws.Cells.Item(40, 8).Value = "IN PROGRESS"

# I40 - Finalized date: still empty, but keep the date/wrap style.
$ws.Cells.Item(40, 9).NumberFormat = "m/d/yy h:mm"
$ws.Cells.Item(40, 9).WrapText = $true

# Row height matches the neighbouring data rows.
$ws.Rows.Item(40).RowHeight = $ws.Rows.Item(39).RowHeight

# ------------------------------------------------------------------
# Window / selection state, mirroring the scrolled/selected cell the
# workbook was left on.
# ------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A45").Select() | Out-Null
